# Neue Rechnung eingefügt Fliegende Klassenzimmer
# Adds a new invoice row ("Das Fliegende Klassenzimmer") to the "Ausgaben"
# table (Table16) on the Ausgaben sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ausgaben")
$ws.Activate()

$tbl = $ws.ListObjects.Item("Table16")

# Adding a ListRow grows the table range (and dimension) automatically,
# matching the A1:I14 -> A1:I15 change in the diff.
$newRow = $tbl.ListRows.Add()
$r = $newRow.Range.Row

# --- Plain text columns -------------------------------------------------
$ws.Cells.Item($r, 1).Value = "Film"
$ws.Cells.Item($r, 2).Value = "Rechnung Film Das Fliegende Klassenzimmer"
$ws.Cells.Item($r, 5).Value = "DCM Film Distribution (Schweiz) GmbH"
$ws.Cells.Item($r, 6).Value = "Kreuzstrasse 2, 8008 Zürich"
$ws.Cells.Item($r, 7).Value = "96 57660 00000 00000 00301 57626"

# --- Datum (column C) - reuse the existing date format from the row above
$ws.Cells.Item(14, 3).Copy()
$ws.Cells.Item($r, 3).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item($r, 3).Value = 45279        # 19.12.2023

# --- Betrag (column D) - reuse the existing CHF currency format
$ws.Cells.Item(14, 4).Copy()
$ws.Cells.Item($r, 4).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item($r, 4).Value = 161.55

# --- Rechnungsnummer (column H) - keep as text so the leading zero stays
$ws.Cells.Item(14, 8).Copy()
$ws.Cells.Item($r, 8).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item($r, 8).Value = "03015762"

# --- Spieldatum (column I) - reuse the existing date format
$ws.Cells.Item(12, 9).Copy()
$ws.Cells.Item($r, 9).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item($r, 9).Value = 45242        # 12.11.2023

$excel.CutCopyMode = 0

# Mirror the author's post-edit selection (cell below/right of the new row).
$ws.Range("C16").Select()
